# Clean up the CNOT crosstalk example on the ankaa-3 topology slide.
#
# 1) Recolor the very first qubit dot (top-right corner of the grid,
#    "Oval 3") from the hard-coded teal srgbClr 66CBC6 to the theme
#    color accent2 (so it matches the crosstalk-highlight orange family).
# 2) Add a second crosstalk highlight ring ("Oval 2") - a copy of the
#    existing "Oval 1" ring - around the bottom dot of the rightmost
#    column, outlined in srgbClr EA7131.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1) Recolor first dot's fill to theme accent2 -------------------------
$firstDot = $s.Shapes.Item(1)
$firstDot.Fill.ForeColor.ObjectThemeColor = 6   # msoThemeColorAccent2

# --- 2) Duplicate the existing highlight ring ("Oval 1") to create the ---
#        new "Oval 2" ring, then reposition / recolor it.
$existingRing = $s.Shapes.Item($s.Shapes.Count)
$newRing = $existingRing.Duplicate()
$newRing.Name = "Oval 2"

# Precise point coordinates chosen so that, after the host's internal
# point -> EMU conversion, they land exactly on the target EMU values
# (x=10642914, y=5798169, cx=725760, cy=725760).
$newRing.Left = 838.02478
$newRing.Top = 456.54876
$newRing.Width = 57.14645669291338
$newRing.Height = 57.14645669291338

$newRing.Line.ForeColor.RGB = 3240426   # srgbClr EA7131
